# Generate Report for Handback
# - Status moves from "Ready for handoff" to "Handed back: in sync with
#   en-US" everywhere it appears: Overview!E2/F2 (zh-cn/de-de columns) and
#   the Status column (C2) on both the "zh-cn" and "de-de" detail sheets.
# - "zh-cn" sheet: Latest Handback DateTime (K2) advances to the new
#   handback timestamp, and the Error Detail (P2) is cleared now that the
#   handback is in sync.
# - "de-de" sheet: same two updates, with its own handback timestamp.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "Handed back: in sync with en-US"
$overview.Range("F2").Value = "Handed back: in sync with en-US"
$overview.Columns("E:F").AutoFit()

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "Handed back: in sync with en-US"
$zhcn.Range("K2").Value = "2016-08-27 14:48:38"
# Clear the error detail but keep it a text cell (empty string), matching
# the other "no value" text cells on this table rather than deleting it.
$zhcn.Range("P2").Value = "'"
$zhcn.Range("P2").Style = "Normal"
$zhcn.Columns("C:C").AutoFit()
$zhcn.Columns("P:P").AutoFit()

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "Handed back: in sync with en-US"
$dede.Range("K2").Value = "2016-08-27 14:48:45"
$dede.Range("P2").Value = "'"
$dede.Range("P2").Style = "Normal"
$dede.Columns("C:C").AutoFit()
$dede.Columns("P:P").AutoFit()
